$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Edit an existing "Chat" entry (column A) in place ---
# Row 235: "usergreet My name is username. Is gservice avaialble"
#       -> "usergreet My name is username. Is gservice avaialble to have"
$ws.Range("A235").Value = "usergreet My name is username. Is gservice avaialble to have"

# --- Append three new rows of reservation-related chat data (Group 5) ---
$ws.Range("A302").Value = "She would like to get gservice"
$ws.Range("A302").HorizontalAlignment = -4131
$ws.Range("B302").Value = 5

$ws.Range("A303").Value = "He need to have "
$ws.Range("A303").HorizontalAlignment = -4131
$ws.Range("B303").Value = 5

$ws.Range("A304").Value = "Could you please add appointment for me"
$ws.Range("A304").HorizontalAlignment = -4131
$ws.Range("B304").Value = 5

# --- Edit another existing "Chat" entry (column A) in place ---
# Row 47: "interjection Good day" -> "interjection Good day to you"
$ws.Range("A47").Value = "interjection Good day to you"

# --- Update the view/selection state to match the author's final position ---
$ws.Range("F304").Select()
